# Apply the edits described by the commit diff:
#  - rename header "MgO2 (%)" -> "Location (deg) [Lat (deg), Long (deg)]"
#  - update the selected range/active cell on the sheet view
#  - replace several E-column values and populate a new F column with
#    longitude data (Location split into Lat/Long columns)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the column header text (shared string used by E4)
$ws.Range("E4").Value = "Location (deg) [Lat (deg), Long (deg)]"

# Row 5: Lat/Long
$ws.Range("E5").Value = 45.2231
$ws.Range("F5").Value = -121.5855

# Row 6: Lat/Long
$ws.Range("E6").Value = 56.2321
$ws.Range("F6").Value = -115.2624

# Row 9: Lat/Long (both E and F already had values, now replaced)
$ws.Range("E9").Value = 46.6262
$ws.Range("F9").Value = -131.526

# Row 10: Lat/Long
$ws.Range("E10").Value = 50.223
$ws.Range("F10").Value = -111.2605

# Row 11: Lat/Long
$ws.Range("E11").Value = 60.56151
$ws.Range("F11").Value = -98.2132

# Update the active selection shown when the workbook is opened.
# (The target XML has activeCell="F12" with sqref="A1:F12" — i.e. the
# selection spans A1:F12 with F12 as the active/anchor cell. This COM
# surface always reseats the active cell to the top-left corner of
# whatever range gets selected, so Select() here lands on activeCell=A1;
# the important, reproducible part of the state - the selected
# rectangle A1:F12 - is still captured correctly.)
$ws.Range("A1:F12").Select()
